$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "pietro" column header
$ws.Range("F1").Value = "pietro"

# Per-row "pietro" values for rows 2..17
$values = @(1, 2, 3, 2, 2, 3, 1, 1, 2, 1, 2, 3, 3, 2, 1, 2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

# Update the active selection to match the recorded edit
$ws.Range("H10").Select()
